$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append row 2 (A2:N2) to the procurement plan table.
# A2 and D2 are numeric; the rest are text. H2 ("True") and the date-like
# strings in I2:N2 must stay literal text rather than being auto-converted
# to a boolean/date by Excel's type inference, so they are entered with a
# leading apostrophe (quote prefix) - the standard Excel mechanism for
# forcing text entry - and the style is reset back to Normal afterward so
# no visible formatting change remains.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "General Staff Salaries"
$ws.Range("C2").Value = "UGX"
$ws.Range("D2").Value = 10000000
$ws.Range("E2").Value = "GOU"
$ws.Range("F2").Value = "Opening Domestic Bidding"
$ws.Range("G2").Value = "Admeasurement"
$ws.Range("H2").Value = "'True"
$ws.Range("I2").Value = "'2020-09-01"
$ws.Range("J2").Value = "'2020-09-21"
$ws.Range("K2").Value = "'2020-09-28"
$ws.Range("L2").Value = "'2020-10-18"
$ws.Range("M2").Value = "'2020-10-23"
$ws.Range("N2").Value = "'2021-06-30"

$ws.Range("H2:N2").Style = "Normal"
